$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Cell values for the new "continuous frames" annotation columns G:K ---
$ws.Range("G1").Value = "0-800"
$ws.Range("H1").Value = "880-18330"
$ws.Range("I1").Value = "18430-23410"
$ws.Range("J1").Value = "23500-end"
$ws.Range("G2").Value = "0-1690"
$ws.Range("H2").Value = "1800-15160"
$ws.Range("I2").Value = "15310-28360"
$ws.Range("J2").Value = "28450-31130"
$ws.Range("K2").Value = "31190-end"
$ws.Range("G3").Value = "0-3190"
$ws.Range("H3").Value = "3290-21980"
$ws.Range("I3").Value = "22100-end"
$ws.Range("G4").Value = "0-5740"
$ws.Range("H4").Value = "5790-end"
$ws.Range("G5").Value = "0-7590"
$ws.Range("H5").Value = "7650-end"
$ws.Range("G6").Value = "0-640"
$ws.Range("H6").Value = "680-12460"
$ws.Range("I6").Value = "12540-end"
$ws.Range("G7").Value = "0-1820"
$ws.Range("H7").Value = "1890-14140"
$ws.Range("I7").Value = "14220-22210"
$ws.Range("J7").Value = "22320-end"
$ws.Range("G8").Value = "0-end"
$ws.Range("G9").Value = "0-23200"
$ws.Range("H9").Value = "23350-end"
$ws.Range("G10").Value = "0-16940"
$ws.Range("H10").Value = "17100-21300"
$ws.Range("I10").Value = "21430-25020"
$ws.Range("J10").Value = "25110-29680"
$ws.Range("K10").Value = "30050-end"
$ws.Range("G11").Value = "0-13770"
$ws.Range("H11").Value = "13910-22370"
$ws.Range("I11").Value = "22510-25640"
$ws.Range("J11").Value = "25770-end"

# --- 2) Font formatting ---
# Build two helper "seed" cells off-sheet, then copy their formats (not values) onto
# the target ranges via PasteSpecial(xlPasteFormats). Seed2 is derived from seed1 so it
# never revisits the default theme font, keeping the style table minimal.
$seed1 = $ws.Range("M1")
$seed1.Font.Name = "Arial"
$seed1.Font.Size = 10

$seed2 = $ws.Range("M2")
$seed1.Copy()
$seed2.PasteSpecial(-4122)
$seed2.Font.Color = 0

$seed1.Copy()
$ws.Range("G1:K7").PasteSpecial(-4122)
$ws.Range("H8:K8").PasteSpecial(-4122)
$ws.Range("I9:K9").PasteSpecial(-4122)
$ws.Range("K11").PasteSpecial(-4122)

$seed2.Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G9:H9").PasteSpecial(-4122)
$ws.Range("G10:K10").PasteSpecial(-4122)
$ws.Range("G11:J11").PasteSpecial(-4122)

$excel.CutCopyMode = 0
$seed1.Clear()
$seed2.Clear()

# --- 3) Selection / view bookkeeping to match the saved workbook state ---
$ws.Range("K11").Select()
